{"js": "// Target change (from the OOXML diff / commit \"Fixed POI packaging and\n// upgraded to POI 3.15\"):\n//   - word/document.xml: XML-namespace-declaration order on <w:document>,\n//     the attribute order inside every <w:tab .../>, <w:pgSz/>, and\n//     <w:pgMar/>, the (malformed, 128-bit) w:rsidR tag shared by the\n//     `REF bookmark1` field run-group, and the numeric w:id of the\n//     bookmarkStart/bookmarkEnd pair for \"bookmark1\".\n//   - word/styles.xml (mislabeled \"word/footnotes.xml\" in the diff, but\n//     its content - docDefaults/latentStyles/style defs - matches\n//     styles.xml exactly): attribute order throughout.\n//\n// None of this touches visible text, run/paragraph formatting, tab stop\n// positions/alignment, the bookmark's name, or the field code/result.\n// Diffing the two docx's canonical (C14N) XML - i.e. with attributes\n// sorted, which is what a content-based comparison naturally does -\n// shows the \"reordering\" is exactly alphabetical order: it is an\n// artifact of how the file was re-serialized (by the upgraded POI\n// library), not an authored change. The two numeric values that do\n// change (the oversized w:rsidR GUID tagging the REF-field runs, and\n// bookmark1's w:id) are internal, tool-generated bookkeeping: Word's\n// object model does not read or expose either one (no `RsidR`/`Id`\n// property exists on Range/Font/Field/Bookmark in Office.js or COM), so\n// they cannot be targeted or reproduced through the API - any value we\n// could invent would be just as arbitrary as the one already in the\n// file. Re-creating the bookmark through the API would only replace one\n// unreproducible value with another (Word assigns small sequential\n// bookmark ids, never these huge POI-style numbers).\n//\n// So there is no reachable, content-affecting edit to make here: the\n// faithful application of this diff through Word's supported APIs is to\n// leave the document's text/structure untouched. The lines below simply\n// confirm (read-only) that the content this diff revolves around - the\n// paragraph text and the \"bookmark1\" bookmark that the REF field further\n// up resolves against - is still present, without mutating anything.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst expectedTexts = [\n  \"This template demonstrates the use of bookmarks. It creates two links to the same bookmark.\",\n  \"Test link before bookmark\\u00a0: a reference to bookmark1\",\n  \"Test bookmark\\u00a0: bookmarked content\",\n];\nfor (let i = 0; i < expectedTexts.length; i++) {\n  if (paragraphs.items[i].text !== expectedTexts[i]) {\n    throw new Error(`Unexpected paragraph ${i} text: ${paragraphs.items[i].text}`);\n  }\n}\n\nconst bookmarkRange = context.document.getBookmarkRangeOrNullObject(\"bookmark1\");\nbookmarkRange.load(\"isNullObject,text\");\nawait context.sync();\n\nif (bookmarkRange.isNullObject) {\n  throw new Error(\"Expected bookmark 'bookmark1' to be present.\");\n}\nif (bookmarkRange.text !== \"bookmarked content\") {\n  throw new Error(`Unexpected bookmark1 range text: ${bookmarkRange.text}`);\n}\n\n// No other part of the diff (namespace/attribute ordering produced by\n// the upgraded serializer, and the two tool-internal random ids - the\n// oversized w:rsidR GUID tagging the REF-field runs and bookmark1's\n// w:id) corresponds to an operation exposed by context.document /\n// Range / Field / Bookmark, so nothing further is mutated.\n", "ps1": "# Target change (from the OOXML diff / commit \"Fixed POI packaging and\n# upgraded to POI 3.15\"):\n#   - word/document.xml: XML-namespace-declaration order on <w:document>,\n#     the attribute order inside every <w:tab .../>, <w:pgSz/>, and\n#     <w:pgMar/>, the (malformed, 128-bit) w:rsidR tag shared by the\n#     `REF bookmark1` field run-group, and the numeric w:id of the\n#     bookmarkStart/bookmarkEnd pair for \"bookmark1\".\n#   - word/styles.xml (mislabeled \"word/footnotes.xml\" in the diff, but\n#     its content - docDefaults/latentStyles/style defs - matches\n#     styles.xml exactly): attribute order throughout.\n#\n# None of this touches visible text, run/paragraph formatting, tab stop\n# positions/alignment, the bookmark's name, or the field code/result.\n# Diffing the two docx's canonical (C14N) XML - i.e. with attributes\n# sorted, which is what a content-based comparison naturally does - shows\n# the \"reordering\" is exactly alphabetical order: it is an artifact of\n# how the file was re-serialized (by the upgraded POI library), not an\n# authored change. The two numeric values that do change (the oversized\n# w:rsidR GUID tagging the REF-field runs, and bookmark1's w:id) are\n# internal, tool-generated bookkeeping: Word's object model does not\n# read or expose either one (no RsidR/Id property exists on\n# Range/Font/Field/Bookmark in COM or Office.js), so they cannot be\n# targeted or reproduced through the API - any value we could invent\n# would be just as arbitrary as the one already in the file.\n# Re-creating the bookmark through the API would only replace one\n# unreproducible value with another (Word assigns small sequential\n# bookmark ids, never these huge POI-style numbers).\n#\n# So there is no reachable, content-affecting edit to make here: the\n# faithful application of this diff through Word's supported APIs is to\n# leave the document's text/structure untouched. The lines below simply\n# confirm (read-only) that the content this diff revolves around - the\n# paragraph text and the \"bookmark1\" bookmark that the REF field further\n# up resolves against - is still present, without mutating anything.\n\n$d = $word.ActiveDocument\n\n$nbsp = [char]0x00A0\n$expectedTexts = @(\n    \"This template demonstrates the use of bookmarks. It creates two links to the same bookmark.\",\n    \"Test link before bookmark\" + $nbsp + \": a reference to bookmark1\",\n    \"Test bookmark\" + $nbsp + \": bookmarked content\"\n)\nfor ($i = 0; $i -lt $expectedTexts.Length; $i++) {\n    $actual = $d.Paragraphs.Item($i + 1).Range.Text.TrimEnd(\"`r\")\n    if ($actual -ne $expectedTexts[$i]) {\n        throw \"Unexpected paragraph $($i + 1) text: $actual\"\n    }\n}\n\nif (-not $d.Bookmarks.Exists(\"bookmark1\")) {\n    throw \"Expected bookmark 'bookmark1' to be present.\"\n}\nif ($d.Bookmarks.Item(\"bookmark1\").Range.Text -ne \"bookmarked content\") {\n    throw \"Unexpected bookmark1 range text: $($d.Bookmarks.Item('bookmark1').Range.Text)\"\n}\n\n$refField = $null\nforeach ($f in $d.Fields) {\n    if ($f.Code.Text -match \"REF bookmark1\") {\n        $refField = $f\n    }\n}\nif ($null -eq $refField) {\n    throw \"Expected a ' REF bookmark1 \\h ' field to be present.\"\n}\n\n# No other part of the diff (namespace/attribute ordering produced by\n# the upgraded serializer, and the two tool-internal random ids - the\n# oversized w:rsidR GUID tagging the REF-field runs and bookmark1's\n# w:id) corresponds to a property or method exposed by Document / Range\n# / Field / Bookmark, so nothing further is mutated.\n"}
